# edit.ps1 - COM-interop script reproducing the authored diff:
#   1. Refresh the cached "datetimeFigureOut" date field text from
#      2/15/26 -> 2/17/26 on the slide master and on every slide layout
#      (this is the auto-date placeholder PowerPoint re-caches on save).
#   2. Update the "problem statement" textbox on slide 2 (shape
#      "TextBox 7"): swap the example store names (PakStyle -> telemart,
#      iShopping -> shopive) while keeping "Daraz" as its own run; the
#      shape's height then auto-fits to the new wrapped text.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# 1) Slide master's Date Placeholder.
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "2/17/26"

# 2) Every slide layout's Date Placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "2/17/26"
}

# 3) Slide 2 - "TextBox 7": update store names in the problem statement.
$slide2 = $p.Slides.Item(2)
$shape = $slide2.Shapes.Item("TextBox 7")
$tr = $shape.TextFrame.TextRange

# Original: "...stores like Daraz, PakStyle, and iShopping."
# Target:   "...stores like Daraz, telemart, and shopive."
# Replace right-to-left so earlier character offsets stay valid, and
# re-assign "Daraz" to itself so it becomes its own run (matching the
# authored run split).
$fullText = $tr.Text
$iShoppingStart = $fullText.IndexOf("iShopping") + 1
$tr.Characters($iShoppingStart, 9).Text = "shopive"

$fullText = $tr.Text
$pakStyleStart = $fullText.IndexOf("PakStyle") + 1
$tr.Characters($pakStyleStart, 8).Text = "telemart"

$fullText = $tr.Text
$darazStart = $fullText.IndexOf("Daraz") + 1
$tr.Characters($darazStart, 5).Text = "Daraz"

Write-Host "Slide2 TextBox7 ->" $tr.Text
